$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "20160406_104403"
$ws.Range("B13").Value = 1196.616
$ws.Range("C13").Value = "remove multiple spaces, convert to lower, convert unicode to ascii, trim `"space`" and `",`""
$ws.Range("D13").Value = "13 features: length, #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2_0, #max_digit_skip_0_2_1, #max_digit_skip_0_2_2, first_character_type_0, first_character_type_1, first_character_type_2, first_character_type_3, #`"space`""
$ws.Range("E13").Value = "Neuron Network"
$ws.Range("F13").Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 0.963696369636964
$ws.Range("I13").Value = "0 filters: "
$ws.Range("J13").Value = 0.0689655172413793

$ws.Range("A14").Value = "20160406_110400"
$ws.Range("B14").Value = 1195.024
$ws.Range("C14").Value = "remove multiple spaces, convert to lower, convert unicode to ascii, trim `"space`" and `",`""
$ws.Range("D14").Value = "13 features: length, #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2_0, #max_digit_skip_0_2_1, #max_digit_skip_0_2_2, first_character_type_0, first_character_type_1, first_character_type_2, first_character_type_3, #`"space`""
$ws.Range("E14").Value = "Neuron Network"
$ws.Range("F14").Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"
$ws.Range("G14").Value = 0.999333333333333
$ws.Range("H14").Value = 0.947194719471947
$ws.Range("I14").Value = "0 filters: "
$ws.Range("J14").Value = 0.0240963855421687

$ws.Range("A15").Value = "20160406_112355"
$ws.Range("B15").Value = 1227.001
$ws.Range("C15").Value = "remove multiple spaces, convert to lower, convert unicode to ascii, trim `"space`" and `",`""
$ws.Range("D15").Value = "13 features: length, #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2_0, #max_digit_skip_0_2_1, #max_digit_skip_0_2_2, first_character_type_0, first_character_type_1, first_character_type_2, first_character_type_3, #`"space`""
$ws.Range("E15").Value = "Neuron Network"
$ws.Range("F15").Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0.96039603960396
$ws.Range("I15").Value = "0 filters: "
$ws.Range("J15").Value = 0.0813953488372093

$ws.Range("A16").Value = "20160406_114422"
$ws.Range("B16").Value = 1250.791
$ws.Range("C16").Value = "remove multiple spaces, convert to lower, convert unicode to ascii, trim `"space`" and `",`""
$ws.Range("D16").Value = "13 features: length, #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2_0, #max_digit_skip_0_2_1, #max_digit_skip_0_2_2, first_character_type_0, first_character_type_1, first_character_type_2, first_character_type_3, #`"space`""
$ws.Range("E16").Value = "Neuron Network"
$ws.Range("F16").Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"
$ws.Range("G16").Value = 0.998666666666667
$ws.Range("H16").Value = 0.95049504950495
$ws.Range("I16").Value = "0 filters: "
$ws.Range("J16").Value = 0.0357142857142857

$ws.Range("A17").Value = "20160406_120513"
$ws.Range("B17").Value = 1289.39
$ws.Range("C17").Value = "remove multiple spaces, convert to lower, convert unicode to ascii, trim `"space`" and `",`""
$ws.Range("D17").Value = "13 features: length, #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2_0, #max_digit_skip_0_2_1, #max_digit_skip_0_2_2, first_character_type_0, first_character_type_1, first_character_type_2, first_character_type_3, #`"space`""
$ws.Range("E17").Value = "Neuron Network"
$ws.Range("F17").Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 1000"
$ws.Range("G17").Value = 0.999333333333333
$ws.Range("H17").Value = 0.953795379537954
$ws.Range("I17").Value = "0 filters: "
$ws.Range("J17").Value = 0.0714285714285714

$ws.Range("A18").Value = "20160406_133507"
$ws.Range("B18").Value = 2757.21
$ws.Range("C18").Value = "trim `"space`" and `",`", convert unicode to ascii, convert to lower, remove multiple spaces"
$ws.Range("D18").Value = "13 features: length, #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2_0, #max_digit_skip_0_2_1, #max_digit_skip_0_2_2, first_character_type_0, first_character_type_1, first_character_type_2, first_character_type_3, #`"space`""
$ws.Range("E18").Value = "Neuron Network"
$ws.Range("F18").Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"
$ws.Range("G18").Value = 0.998666666666667
$ws.Range("H18").Value = 0.957095709570957
$ws.Range("I18").Value = "0 filters: "
$ws.Range("J18").Value = 0.0588235294117647

$ws.Range("A19").Value = "20160406_142105"
$ws.Range("B19").Value = 2768.523
$ws.Range("C19").Value = "trim `"space`" and `",`", convert unicode to ascii, convert to lower, remove multiple spaces"
$ws.Range("D19").Value = "13 features: length, #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2_0, #max_digit_skip_0_2_1, #max_digit_skip_0_2_2, first_character_type_0, first_character_type_1, first_character_type_2, first_character_type_3, #`"space`""
$ws.Range("E19").Value = "Neuron Network"
$ws.Range("F19").Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"
$ws.Range("G19").Value = 0.999333333333333
$ws.Range("H19").Value = 0.963696369636964
$ws.Range("I19").Value = "0 filters: "
$ws.Range("J19").Value = 0.0689655172413793

$ws.Range("A20").Value = "20160406_150713"
$ws.Range("B20").Value = 1689.47
$ws.Range("C20").Value = "trim `"space`" and `",`", convert unicode to ascii, convert to lower, remove multiple spaces"
$ws.Range("D20").Value = "13 features: length, #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2_0, #max_digit_skip_0_2_1, #max_digit_skip_0_2_2, first_character_type_0, first_character_type_1, first_character_type_2, first_character_type_3, #`"space`""
$ws.Range("E20").Value = "Neuron Network"
$ws.Range("F20").Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"
$ws.Range("G20").Value = 0.998666666666667
$ws.Range("H20").Value = 0.947194719471947
$ws.Range("I20").Value = "0 filters: "
$ws.Range("J20").Value = 0.036144578313253

$ws.Range("A21").Value = "20160406_153523"
$ws.Range("B21").Value = 1622.853
$ws.Range("C21").Value = "trim `"space`" and `",`", convert unicode to ascii, convert to lower, remove multiple spaces"
$ws.Range("D21").Value = "13 features: length, #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2_0, #max_digit_skip_0_2_1, #max_digit_skip_0_2_2, first_character_type_0, first_character_type_1, first_character_type_2, first_character_type_3, #`"space`""
$ws.Range("E21").Value = "Neuron Network"
$ws.Range("F21").Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"
$ws.Range("G21").Value = 0.999333333333333
$ws.Range("H21").Value = 0.95049504950495
$ws.Range("I21").Value = "0 filters: "
$ws.Range("J21").Value = 0

$ws.Range("A22").Value = "20160406_160226"
$ws.Range("B22").Value = 1509.281
$ws.Range("C22").Value = "trim `"space`" and `",`", convert unicode to ascii, convert to lower, remove multiple spaces"
$ws.Range("D22").Value = "13 features: length, #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_2_0, #max_digit_skip_0_2_1, #max_digit_skip_0_2_2, first_character_type_0, first_character_type_1, first_character_type_2, first_character_type_3, #`"space`""
$ws.Range("E22").Value = "Neuron Network"
$ws.Range("F22").Value = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0.95049504950495
$ws.Range("I22").Value = "0 filters: "
$ws.Range("J22").Value = 0.0833333333333333

